# modifs to smart random
#
# Applies to slide 38 ("Collected data" / smart-random explanation slide):
#   - grows & edits the STATES box (TextBox 2) - adds an "Exit_remaining" line
#     and extends the "Taken_transitions" line
#   - slides the DEADLOCKS box (TextBox 6) to the right to make room
#   - moves/grows the DEAD_PREDS/TR_COUNT/TR_COUPLES box (TextBox 7) and moves
#     the "=> % of each transition overall" phrase up onto the TR_COUNT line
# Also marks the final slide (41) as hidden from the slide show.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(38)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# TextBox 2 ("STATES" legend box): add the "Exit_remaining" line, extend the
# "Taken_transitions" line, and grow the box to fit the extra text.
# ---------------------------------------------------------------------------
$statesBox = Get-ShapeByName $s "TextBox 2"

$newLine = [char]13
$statesText = "STATES" + $newLine +
    "State : explicit state" + $newLine +
    "Seen : how many times it was visited" + $newLine +
    "Exit_number : how many exits at this state" + $newLine +
    "Exit_transitions : which exits possible" + $newLine +
    "Exit_remaining : how many left to take" + $newLine +
    "Taken_transitions : which exits were taken and how often"

$statesBox.TextFrame2.TextRange.Text = $statesText

$statesBox.Width = 497.53641732283467
$statesBox.Height = 142.98287401574802

# ---------------------------------------------------------------------------
# TextBox 6 ("DEADLOCKS" box): shift right to make room for the wider
# STATES box (text content is unchanged).
# ---------------------------------------------------------------------------
$deadlocksBox = Get-ShapeByName $s "TextBox 6"

$deadlocksBox.Left = 414.2661811023622
$deadlocksBox.Top = 288.7526377952756

# ---------------------------------------------------------------------------
# TextBox 7 (DEAD_PREDS / TR_COUNT / TR_COUPLES box): move the "=> % of each
# transition overall" phrase from the TR_COUPLES line up onto the TR_COUNT
# line, move the box down and widen it.
# ---------------------------------------------------------------------------
$trCountBox = Get-ShapeByName $s "TextBox 7"

$trCountText = "DEAD_PREDS: state -> state [states that had a deadlock follow them]" + $newLine +
    "TR_COUNT: transitions that appeared [and how many times] => % of each transition overall " + $newLine +
    "TR_COUPLES : how often ti -> tj appears" + $newLine +
    $newLine

$trCountBox.TextFrame2.TextRange.Text = $trCountText

$trCountBox.Left = 7.436496062992126
$trCountBox.Top = 444.4297244094488
$trCountBox.Width = 636.6716141732284
$trCountBox.Height = 84.82035433070867

# ---------------------------------------------------------------------------
# Hide the last slide (41) from the slide show.
# ---------------------------------------------------------------------------
$lastSlide = $p.Slides.Item(41)
$lastSlide.SlideShowTransition.Hidden = $true
